$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 117, shifting rows 117:183 down to 118:184
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new data point
$ws.Cells.Item(117, 1).Value = 1
$ws.Cells.Item(117, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(117, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(117, 4).Value = 44574
$ws.Cells.Item(117, 5).Value = 15
$ws.Cells.Item(117, 6).Value = "Fruta"
$ws.Cells.Item(117, 7).Value = 100102
$ws.Cells.Item(117, 8).Value = "Cítricos"
$ws.Cells.Item(117, 9).Value = 100102003
$ws.Cells.Item(117, 10).Value = "Limón"
$ws.Cells.Item(117, 11).Value = "Sin especificar"
$ws.Cells.Item(117, 12).Value = "2a amarillo"
$ws.Cells.Item(117, 13).Value = 250
$ws.Cells.Item(117, 14).Value = 25000
$ws.Cells.Item(117, 15).Value = 26000
$ws.Cells.Item(117, 16).Value = 25500
$ws.Cells.Item(117, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(117, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(117, 19).Value = 1275
$ws.Cells.Item(117, 20).Value = 20
